$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the kernel/padding values on row 3 (C3: 8 -> 4, D3: 3 -> 1)
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 1

# Update the padding/stride values on row 7 (D7: 1 -> 0, E7: 2 -> 1)
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 1

# Delete row 8 entirely (it is no longer needed now that F7 rounds down to 1)
$ws.Rows("8:8").Delete()

# The column-A width nudges slightly (content reflow after the row delete)
$ws.Columns("A:A").ColumnWidth = 9.25

# Reflect the final selection recorded in the saved file
$ws.Range("B3").Select()
